# This script reproduces the diff: the existing 2 match rows are replaced by
# an updated set of 11 match rows (rows 2-12) of Rahul Tripathi IPL 2020 stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new blank rows starting at row 2, pushing the 2 existing data rows
# down out of the way; we will overwrite every data row (2-12) below anyway.
$ws.Rows("2:10").Insert()

# Row 2
$ws.Range("A2").Value = " Oct 21 2020"
$ws.Range("B2").Value = " Abu Dhabi"
$ws.Range("C2").Value = "RCB won by 8 wickets (with 39 balls remaining)"
$ws.Range("D2").Value = "Kolkata Knight Riders"
$ws.Range("E2").Value = "Royal Challengers Bangalore"
$ws.Range("F2").Value = "Rahul Tripathi "
$ws.Range("G2").Value = "'1"
$ws.Range("H2").Value = "'5"
$ws.Range("I2").Value = "'0"
$ws.Range("J2").Value = "'0"
$ws.Range("K2").Value = "'20.00"

# Row 3
$ws.Range("A3").Value = " Oct 12 2020"
$ws.Range("B3").Value = " Sharjah"
$ws.Range("C3").Value = "RCB won by 82 runs"
$ws.Range("D3").Value = "Kolkata Knight Riders"
$ws.Range("E3").Value = "Royal Challengers Bangalore"
$ws.Range("F3").Value = "Rahul Tripathi "
$ws.Range("G3").Value = "'16"
$ws.Range("H3").Value = "'22"
$ws.Range("I3").Value = "'1"
$ws.Range("J3").Value = "'0"
$ws.Range("K3").Value = "'72.72"

# Row 4
$ws.Range("A4").Value = " Oct 16 2020"
$ws.Range("B4").Value = " Abu Dhabi"
$ws.Range("C4").Value = "Mumbai won by 8 wickets (with 19 balls remaining)"
$ws.Range("D4").Value = "Kolkata Knight Riders"
$ws.Range("E4").Value = "Mumbai Indians"
$ws.Range("F4").Value = "Rahul Tripathi "
$ws.Range("G4").Value = "'7"
$ws.Range("H4").Value = "'9"
$ws.Range("I4").Value = "'1"
$ws.Range("J4").Value = "'0"
$ws.Range("K4").Value = "'77.77"

# Row 5
$ws.Range("A5").Value = " Oct 3 2020"
$ws.Range("B5").Value = " Sharjah"
$ws.Range("C5").Value = "Capitals won by 18 runs"
$ws.Range("D5").Value = "Kolkata Knight Riders"
$ws.Range("E5").Value = "Delhi Capitals"
$ws.Range("F5").Value = "Rahul Tripathi "
$ws.Range("G5").Value = "'36"
$ws.Range("H5").Value = "'16"
$ws.Range("I5").Value = "'3"
$ws.Range("J5").Value = "'3"
$ws.Range("K5").Value = "'225.00"

# Row 6
$ws.Range("A6").Value = " Oct 24 2020"
$ws.Range("B6").Value = " Abu Dhabi"
$ws.Range("C6").Value = "KKR won by 59 runs"
$ws.Range("D6").Value = "Kolkata Knight Riders"
$ws.Range("E6").Value = "Delhi Capitals"
$ws.Range("F6").Value = "Rahul Tripathi "
$ws.Range("G6").Value = "'13"
$ws.Range("H6").Value = "'12"
$ws.Range("I6").Value = "'1"
$ws.Range("J6").Value = "'0"
$ws.Range("K6").Value = "'108.33"

# Row 7
$ws.Range("A7").Value = " Oct 18 2020"
$ws.Range("B7").Value = " Abu Dhabi"
$ws.Range("C7").Value = "Match tied (KKR won the one-over eliminator)"
$ws.Range("D7").Value = "Kolkata Knight Riders"
$ws.Range("E7").Value = "Sunrisers Hyderabad"
$ws.Range("F7").Value = "Rahul Tripathi "
$ws.Range("G7").Value = "'23"
$ws.Range("H7").Value = "'16"
$ws.Range("I7").Value = "'2"
$ws.Range("J7").Value = "'1"
$ws.Range("K7").Value = "'143.75"

# Row 8
$ws.Range("A8").Value = " Oct 7 2020"
$ws.Range("B8").Value = " Abu Dhabi"
$ws.Range("C8").Value = "KKR won by 10 runs"
$ws.Range("D8").Value = "Kolkata Knight Riders"
$ws.Range("E8").Value = "Chennai Super Kings"
$ws.Range("F8").Value = "Rahul Tripathi "
$ws.Range("G8").Value = "'81"
$ws.Range("H8").Value = "'51"
$ws.Range("I8").Value = "'8"
$ws.Range("J8").Value = "'3"
$ws.Range("K8").Value = "'158.82"

# Row 9
$ws.Range("A9").Value = " Oct 29 2020"
$ws.Range("B9").Value = " Dubai (DSC)"
$ws.Range("C9").Value = "Super Kings won by 6 wickets"
$ws.Range("D9").Value = "Kolkata Knight Riders"
$ws.Range("E9").Value = "Chennai Super Kings"
$ws.Range("F9").Value = "Rahul Tripathi "
$ws.Range("G9").Value = "'3"
$ws.Range("H9").Value = "'2"
$ws.Range("I9").Value = "'0"
$ws.Range("J9").Value = "'0"
$ws.Range("K9").Value = "'150.00"

# Row 10
$ws.Range("A10").Value = " Oct 26 2020"
$ws.Range("B10").Value = " Sharjah"
$ws.Range("C10").Value = "Kings XI won by 8 wickets (with 7 balls remaining)"
$ws.Range("D10").Value = "Kolkata Knight Riders"
$ws.Range("E10").Value = "Kings XI Punjab"
$ws.Range("F10").Value = "Rahul Tripathi "
$ws.Range("G10").Value = "'7"
$ws.Range("H10").Value = "'4"
$ws.Range("I10").Value = "'0"
$ws.Range("J10").Value = "'1"
$ws.Range("K10").Value = "'175.00"

# Row 11
$ws.Range("A11").Value = " Oct 10 2020"
$ws.Range("B11").Value = " Abu Dhabi"
$ws.Range("C11").Value = "KKR won by 2 runs"
$ws.Range("D11").Value = "Kolkata Knight Riders"
$ws.Range("E11").Value = "Kings XI Punjab"
$ws.Range("F11").Value = "Rahul Tripathi "
$ws.Range("G11").Value = "'4"
$ws.Range("H11").Value = "'10"
$ws.Range("I11").Value = "'1"
$ws.Range("J11").Value = "'0"
$ws.Range("K11").Value = "'40.00"

# Row 12
$ws.Range("A12").Value = " Nov 1 2020"
$ws.Range("B12").Value = " Dubai (DSC)"
$ws.Range("C12").Value = "KKR won by 60 runs"
$ws.Range("D12").Value = "Kolkata Knight Riders"
$ws.Range("E12").Value = "Rajasthan Royals"
$ws.Range("F12").Value = "Rahul Tripathi "
$ws.Range("G12").Value = "'39"
$ws.Range("H12").Value = "'34"
$ws.Range("I12").Value = "'4"
$ws.Range("J12").Value = "'2"
$ws.Range("K12").Value = "'114.70"
